# Fix Training Data Issue (#48)
# The BF column ("Date") held values in the wrong format (e.g. "5-25-2013-14",
# an artifact of how the NBA stats site showed the date) - normalize them to
# ISO "YYYY-MM-DD" (e.g. "2014-05-25").
#
# Excel's automatic "looks like a date" typing would otherwise convert a
# literal "2014-05-25" assignment into a date serial number, so the target
# cells are pre-formatted as Text ("@") before the value is written - this
# keeps them as plain strings, matching the original data's text values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$firstRow = 2
$lastRow = 31
$col = "BF"

$oldValue = "5-25-2013-14"
$newValue = "2014-05-25"

$targetRange = $ws.Range("$col$firstRow`:$col$lastRow")
$targetRange.NumberFormat = "@"

for ($row = $firstRow; $row -le $lastRow; $row++) {
    $cell = $ws.Range("$col$row")
    if ($cell.Value2 -eq $oldValue) {
        $cell.Value = $newValue
    }
}
